$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01514828764759746
$ws.Range("C2").Value = 0.3127903958511391
$ws.Range("D2").Value = 26.21740644021617
$ws.Range("E2").Value = 616238.5361209477
$ws.Range("G2").Value = 616265.0814660714

$ws.Range("B3").Value = 0.04763786555579896
$ws.Range("C3").Value = 0.04240448674262143
$ws.Range("D3").Value = 26.21740644021617
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 34.96768127846357

$ws.Range("B4").Value = 0.3048080303191223
$ws.Range("C4").Value = 0.3127903958511391
$ws.Range("D4").Value = 0.1575252929769615
$ws.Range("E4").Value = 8.660232485948974
$ws.Range("G4").Value = 9.435356205096197

$ws.Range("B5").Value = 3.230985683306322
$ws.Range("C5").Value = 1.667794583268128
$ws.Range("D5").Value = 3.900430680208489
$ws.Range("E5").Value = 8.660232485948974
$ws.Range("G5").Value = 17.45944343273191
